# Updated cryptos list on Sat Sep  7 11:26:55 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for the crypto table,
# and swaps the BitcoinCash / Uniswap rows (21 and 22) to reflect the new
# ranking order, with their refreshed price/volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.449.67'
$ws.Range('D3').Value = '2.286.98'
$ws.Range('E3').Value = '  -2.84%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '493.93'
$ws.Range('E5').Value = '  -2.09%  '
$ws.Range('D6').Value = '127.26'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').Value = '0.528'
$ws.Range('E8').Value = '  -1.69%  '
$ws.Range('D9').Value = '2.284.99'
$ws.Range('E9').Value = '  -3.50%  '
$ws.Range('D10').Value = '0.0944'
$ws.Range('E10').Value = '  -2.98%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('E13').Value = '  -3.87%  '
$ws.Range('D14').Value = '2.690.90'
$ws.Range('E14').Value = '  -2.94%  '
$ws.Range('D15').Value = '21.57'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '54.367.31'
$ws.Range('E16').Value = '  -2.76%  '
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').Value = '2.276.97'
$ws.Range('E18').Value = '  -4.72%  '
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '303.55'
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '6.50'
$ws.Range('E22').Value = '  +4.92%  '
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('E24').Value = '  -2.63%  '
$ws.Range('D25').Value = '63.57'
$ws.Range('E25').Value = '  -2.82%  '
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').Value = '2.391.40'
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('D31').Value = '169.03'
$ws.Range('E31').Value = '  -1.38%  '
$ws.Range('D32').Value = '1.60'
$ws.Range('E32').Value = '  -2.33%  '
$ws.Range('E33').Value = '  -2.90%  '
$ws.Range('D34').Value = '5.88'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('D40').Value = '0.873'
$ws.Range('E40').Value = '  +4.10%  '
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').Value = '35.55'
$ws.Range('E43').Value = '  +1.05%  '
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').Value = '128.24'
$ws.Range('E46').Value = '  +2.12%  '
$ws.Range('D47').Value = '4.81'
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D48').Value = '0.0891'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('D50').Value = '239.55'
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('D51').Value = '0.0479'
$ws.Range('E51').Value = '  +0.24%  '
